# Script 1 - atualizacao em 2025-09-06 02:14:47Z
#
# The dataset drops the "01/01/2014" row for each region (Brasil, Nordeste,
# Sergipe) and shifts the remaining rows up by one, while also refreshing a
# handful of the recomputed rate values for 2016/2020-2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first (01/01/2014) data row of each of the three 11-row blocks.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("24:24").Delete()   # Sergipe - 01/01/2014
$ws.Rows("13:13").Delete()   # Nordeste - 01/01/2014
$ws.Rows("2:2").Delete()     # Brasil - 01/01/2014

# After the deletions, update the recomputed values that differ from the
# simple shift-up of the remaining rows.
$ws.Range("C3").Value  = 0.1997031955321377    # Brasil  01/01/2016
$ws.Range("C7").Value  = 0.3175325379261067    # Brasil  01/01/2020
$ws.Range("C8").Value  = 0.3892553809515489    # Brasil  01/01/2021
$ws.Range("C9").Value  = 0.480189775225401     # Brasil  01/01/2022
$ws.Range("C10").Value = 0.7044517170836323    # Brasil  01/01/2023
$ws.Range("C11").Value = 0.740784122293602     # Brasil  01/01/2024

$ws.Range("C17").Value = 0.468637543760702     # Nordeste 01/01/2020
$ws.Range("C18").Value = 0.5295549002286808    # Nordeste 01/01/2021
$ws.Range("C19").Value = 0.6729915024147284    # Nordeste 01/01/2022
$ws.Range("C20").Value = 0.8970330678180559    # Nordeste 01/01/2023
$ws.Range("C21").Value = 1.008947005262028     # Nordeste 01/01/2024

$ws.Range("C27").Value = 0.4713350981788455    # Sergipe  01/01/2020
$ws.Range("C28").Value = 0.5924473086143022    # Sergipe  01/01/2021
$ws.Range("C29").Value = 0.7948040443998168    # Sergipe  01/01/2022
$ws.Range("C30").Value = 1.115440423224473     # Sergipe  01/01/2023
$ws.Range("C31").Value = 1.197315040200634     # Sergipe  01/01/2024
